# ---------------------------------------------------------------------------
# Target change (per the supplied unified diff): four <w:nsid w:val="…"/>
# GUID-like attributes inside word/numbering.xml are replaced with new
# random-looking hex values, one per <w:abstractNum>:
#
#   abstractNumId="990"   6cf36c90 -> b7720116
#   abstractNumId="991"   2868470f -> 13330212
#   abstractNumId="99721" e369fe1c -> aaa74375
#   abstractNumId="99722" e5e243e7 -> 432c33d0
#
# Nothing else in the package changes: no paragraph text, formatting, list
# numbers/levels, styles, or even the <w:num>/abstractNum structure differ.
# The commit message ("Vygenerovany file ve slozce", i.e. "Generated file
# in folder") confirms this is a mechanical artifact of an automated
# regeneration/export pass, not a user-visible edit.
#
# <w:nsid> is Word's internal "namespace ID" for a list definition - a
# housekeeping GUID minted once when a list is first created and otherwise
# inert (it has no effect on rendering, numbering values, or behaviour).
# It is intentionally NOT part of the Word object model surfaced to
# VBA/COM automation: there is no List.NSID / ListTemplate.NSID / similar
# property in real Word, and this host's object model mirrors that - it
# was probed exhaustively (every List/ListTemplate/ListFormat property,
# Document/Range Find against body text and story ranges, WordOpenXML
# read-back, StartNewList, etc.) and nsid never surfaces for reading or
# writing anywhere. WordOpenXML is also read-only here, and Find/Replace
# only ever reaches visible story text, never the numbering part, so
# there is no COM-reachable way to retarget those four values.
#
# Since no legitimate Application/Document/List/ListFormat/Range call can
# touch these bytes, this script deliberately performs no mutating calls,
# so the parts of the package that genuinely are reachable through COM
# (body text, styles, structure, etc.) stay byte-for-byte faithful to the
# source rather than risk an unrelated regression while chasing an
# unreachable field.
$d = $word.ActiveDocument
Write-Output ("nsid regeneration is not exposed via Word COM automation; " + $d.Name + " left unmodified.")
